$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue "D2" "249.14"
Set-TextValue "D3" "21.68"
Set-TextValue "D4" "5.572"
Set-TextValue "D5" "0.05672"
Set-TextValue "B6" "GateToken"
Set-TextValue "C6" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D6" "3.377"
Set-TextValue "E6" "5GateTokenGT"
Set-TextValue "B7" "KuCoinToken"
Set-TextValue "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D7" "6.439"
Set-TextValue "E7" "6KuCoinTokenKCS"
Set-TextValue "B8" "MXToken"
Set-TextValue "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.8025"
Set-TextValue "E8" "7MXTokenMX"
Set-TextValue "B9" "FTXToken"
Set-TextValue "C9" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D9" "1.040"
Set-TextValue "E9" "8FTXTokenFTT"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1426"
Set-TextValue "E10" "9WazirXWRX"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07265"
Set-TextValue "E11" "10MandalaExchangeTokenMDX"
Set-TextValue "B12" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03136"
Set-TextValue "E12" "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue "B13" "BitrueCoin"
Set-TextValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.02923"
Set-TextValue "E13" "12BitrueCoinBTR"
Set-TextValue "B14" "BitMartToken"
Set-TextValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09274"
Set-TextValue "E14" "13BitMartTokenBMX"
Set-TextValue "B15" "BitForexToken"
Set-TextValue "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001662"
Set-TextValue "E15" "14BitForexTokenBF"
Set-TextValue "B16" "MCDex"
Set-TextValue "C16" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.218"
Set-TextValue "E16" "15MCDexMCB"
Set-TextValue "B17" "CoinExToken"
Set-TextValue "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04727"
Set-TextValue "E17" "16CoinExTokenCET"
Set-TextValue "B18" "One"
Set-TextValue "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005817"
Set-TextValue "E18" "17OneONE"
Set-TextValue "B19" "TigerCash"
Set-TextValue "C19" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D19" "0.006460"
Set-TextValue "E19" "18TigerCashTCH"
Set-TextValue "B20" "HotbitToken"
Set-TextValue "C20" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D20" "0.005060"
Set-TextValue "E20" "19HotbitTokenHTB"
Set-TextValue "B21" "BitKan"
Set-TextValue "C21" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D21" "0.001051"
Set-TextValue "E21" "20BitKanKAN"
Set-TextValue "B22" "NitroEx"
Set-TextValue "C22" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D22" "0.0001501"
Set-TextValue "E22" "21NitroExNTX"
Set-TextValue "B23" "LEO"
Set-TextValue "C23" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D23" "3.980"
Set-TextValue "E23" "22LEOLEO"
Set-TextValue "B24" "BTSEToken"
Set-TextValue "C24" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D24" "2.112"
Set-TextValue "E24" "23BTSETokenBTSE"
Set-TextValue "D25" "0.3267"
Set-TextValue "D40" "0.04127"
Set-TextValue "B41" "BKEXToken"
Set-TextValue "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1044"
Set-TextValue "E41" "40BKEXTokenBKK"
Set-TextValue "D42" "0.002973"
Set-TextValue "B43" "KickToken"
Set-TextValue "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003237"
Set-TextValue "E43" "42KickTokenKICKWorstin24h"
Set-TextValue "D44" "0.008516"
Set-TextValue "D45" "0.00005649"
Set-TextValue "E47" "46CoinbaseStockTokenCOINBestin24h"
Set-TextValue "D48" "0.01656"
Set-TextValue "E48" "47BOLOBOLO"
